# Update "Latest HO Xliff Generate Date" / "Correspond Handoff/Handback Datetime"
# timestamps to reflect a freshly generated handback report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-22 12:47:22"

$wsZhCn.Range("H4").Value = "2016-08-22 12:47:17"
$wsZhCn.Range("K4").Value = "2016-08-22 12:47:36"

$wsDeDe.Range("H4").Value = "2016-08-22 12:47:22"
$wsDeDe.Range("K4").Value = "2016-08-22 12:47:43"
